$wb = $excel.ActiveWorkbook

# --- Sheet "Excel format": update raw date serials (C:F), rows 8-22 ---
$s3 = $wb.Worksheets.Item("Excel format")
$s3.Range("C8:F8").Value = 45233
$s3.Range("C9:F9").Value = 45313
$s3.Range("C10:F10").Value = 45393
$s3.Range("C11:F11").Value = 45410
$s3.Range("C12:F12").Value = 45440
$s3.Range("C13:F13").Value = 45473
$s3.Range("C14:F14").Value = 45553
$s3.Range("C15:F15").Value = 45633
$s3.Range("C16:F16").Value = 45233
$s3.Range("C17:F17").Value = 45313
$s3.Range("C18:F18").Value = 45393
$s3.Range("C19:F19").Value = 45410
$s3.Range("C20:F20").Value = 45440
$s3.Range("C21:F21").Value = 45473
$s3.Range("C22:F22").Value = 45233

function Set-TextDate($rng, $text) {
    $fmt = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = $fmt
}

# --- Sheet "Call metod": local format (F) + German format (G), rows 8-22 ---
$s1 = $wb.Worksheets.Item("Call metod")
Set-TextDate $s1.Range("F8") "2023-November-03"
Set-TextDate $s1.Range("G8") "2023-November-03"
Set-TextDate $s1.Range("F9") "2024-January-22"
Set-TextDate $s1.Range("G9") "2024-Januar-22"
Set-TextDate $s1.Range("F10") "2024-April-11"
Set-TextDate $s1.Range("G10") "2024-April-11"
Set-TextDate $s1.Range("F11") "2024-April-28"
Set-TextDate $s1.Range("G11") "2024-April-28"
Set-TextDate $s1.Range("F12") "2024-May-28"
Set-TextDate $s1.Range("G12") "2024-Mai-28"
Set-TextDate $s1.Range("F13") "2024-June-30"
Set-TextDate $s1.Range("G13") "2024-Juni-30"
Set-TextDate $s1.Range("F14") "2024-September-18"
Set-TextDate $s1.Range("G14") "2024-September-18"
Set-TextDate $s1.Range("F15") "2024-December-07"
Set-TextDate $s1.Range("G15") "2024-Dezember-07"
Set-TextDate $s1.Range("F16") "2023-November-03"
Set-TextDate $s1.Range("G16") "2023-November-03"
Set-TextDate $s1.Range("F17") "2024-January-22"
Set-TextDate $s1.Range("G17") "2024-Januar-22"
Set-TextDate $s1.Range("F18") "2024-April-11"
Set-TextDate $s1.Range("G18") "2024-April-11"
Set-TextDate $s1.Range("F19") "2024-April-28"
Set-TextDate $s1.Range("G19") "2024-April-28"
Set-TextDate $s1.Range("F20") "2024-May-28"
Set-TextDate $s1.Range("G20") "2024-Mai-28"
Set-TextDate $s1.Range("F21") "2024-June-30"
Set-TextDate $s1.Range("G21") "2024-Juni-30"
Set-TextDate $s1.Range("F22") "2023-November-03"
Set-TextDate $s1.Range("G22") "2023-November-03"

# --- Sheet "Alternatives ;cond=": Russian format (F) + US format (G), rows 8-22 ---
$s2 = $wb.Worksheets.Item("Alternatives ;cond=")
Set-TextDate $s2.Range("F8") "03.11.2023"
Set-TextDate $s2.Range("G8") "11/03/2023"
Set-TextDate $s2.Range("F9") "22.01.2024"
Set-TextDate $s2.Range("G9") "01/22/2024"
Set-TextDate $s2.Range("F10") "11.04.2024"
Set-TextDate $s2.Range("G10") "04/11/2024"
Set-TextDate $s2.Range("F11") "28.04.2024"
Set-TextDate $s2.Range("G11") "04/28/2024"
Set-TextDate $s2.Range("F12") "28.05.2024"
Set-TextDate $s2.Range("G12") "05/28/2024"
Set-TextDate $s2.Range("F13") "30.06.2024"
Set-TextDate $s2.Range("G13") "06/30/2024"
Set-TextDate $s2.Range("F14") "18.09.2024"
Set-TextDate $s2.Range("G14") "09/18/2024"
Set-TextDate $s2.Range("F15") "07.12.2024"
Set-TextDate $s2.Range("G15") "12/07/2024"
Set-TextDate $s2.Range("F16") "03.11.2023"
Set-TextDate $s2.Range("G16") "11/03/2023"
Set-TextDate $s2.Range("F17") "22.01.2024"
Set-TextDate $s2.Range("G17") "01/22/2024"
Set-TextDate $s2.Range("F18") "11.04.2024"
Set-TextDate $s2.Range("G18") "04/11/2024"
Set-TextDate $s2.Range("F19") "28.04.2024"
Set-TextDate $s2.Range("G19") "04/28/2024"
Set-TextDate $s2.Range("F20") "28.05.2024"
Set-TextDate $s2.Range("G20") "05/28/2024"
Set-TextDate $s2.Range("F21") "30.06.2024"
Set-TextDate $s2.Range("G21") "06/30/2024"
Set-TextDate $s2.Range("F22") "03.11.2023"
Set-TextDate $s2.Range("G22") "11/03/2023"
